$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2240798473358154
$ws.Range("B3").Value = 0.5015957355499268
$ws.Range("B4").Value = 1.568688631057739
$ws.Range("B5").Value = 3.871842622756958
